$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.890.31"
$ws.Range("E2").Value = "  +0.18%  "

$ws.Range("D3").Value = "2.531.55"
$ws.Range("E3").Value = "  -0.67%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.574"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.83%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0811"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.23%  "

$ws.Range("E13").Value = "  -4.19%  "

$ws.Range("D14").Value = "2.923.99"
$ws.Range("E14").Value = "  -0.51%  "

$ws.Range("D15").Value = "2.486.21"
$ws.Range("E15").Value = "  -2.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.92%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.847"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.18%  "

$ws.Range("D18").Value = "42.963.26"
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.86%  "

$ws.Range("D21").Value = "0.0₃0964"
$ws.Range("E21").Value = "  -0.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.22%  "

$ws.Range("E24").Value = "  -1.01%  "

$ws.Range("E25").Value = "  +0.81%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.06%  "

$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.34%  "

$ws.Range("E29").Value = "  +1.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.35"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.62%  "

$ws.Range("E35").Value = "  +2.84%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0795"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.113"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.23%  "

$ws.Range("E39").Value = "  +1.92%  "

$ws.Range("E40").Value = "  -1.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "21.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.37%  "

$ws.Range("E42").Value = "  -1.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0304"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.27%  "

$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.52%  "

$ws.Range("D46").Value = "1.995.65"
$ws.Range("E46").Value = "  -1.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.06"
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.61%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.779.64"
$ws.Range("E49").Value = "  -0.45%  "

$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.56%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "104.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.65%  "
